# LAB 18 + ConvertExcel2Json
# Update the credentials table on Sheet1: convert the existing numeric
# passwords to text, and append two new username/password rows (with
# "mailto:" hyperlinks on the username cells), matching the new
# ConvertExcel2Json fixture data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- existing rows: passwords become text values -----------------------
$ws.Range("B2").Value = "a12345678"
$ws.Range("B3").Value = "876543s21"

# --- new row 4 -----------------------------------------------------------
$ws.Range("A4").Value = "khoa01@gg.vn"
$ws.Range("B4").Value = "abc123413"

# --- new row 5 -----------------------------------------------------------
$ws.Range("A5").Value = "khoa99@gj.sd"
$ws.Range("B5").Value = "13jfsfsj323"

# --- hyperlink the new username cells (mirrors A2/A3) --------------------
$ws.Hyperlinks.Add($ws.Range("A4"), "mailto:khoa01@gg.vn")
$ws.Hyperlinks.Add($ws.Range("A5"), "mailto:khoa99@gj.sd")

# Hyperlinks.Add re-stamps its own "Hyperlink" cell style; restore the
# same style already used by A2/A3 so the style table doesn't fork.
$ws.Range("A4").Style = $ws.Range("A2").Style
$ws.Range("A5").Style = $ws.Range("A2").Style

# --- selection cursor, as left by the author after editing ---------------
[void]$ws.Range("B10").Select()
